$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 750
$ws.Range("J3").Value = 842
$ws.Range("J4").Value = 182
$ws.Range("J6").Value = 1209
$ws.Range("J7").Value = 3042

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J2").Value = 12
$ws.Range("J4").Value = 4
$ws.Range("J7").Value = 38

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J2").Value = 28
$ws.Range("J3").Value = 35
$ws.Range("J7").Value = 109

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J3").Value = 13
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 40

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 22
$ws.Range("J3").Value = 50
$ws.Range("J7").Value = 116

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J3").Value = 8
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J6").Value = 31
$ws.Range("J7").Value = 77

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J2").Value = 22
$ws.Range("J4").Value = 16
$ws.Range("J6").Value = 28
$ws.Range("J7").Value = 85
$ws.Range("J8").Value = 191
$ws.Range("J9").Value = 22
$ws.Range("J10").Value = 19
$ws.Range("J11").Value = 41
$ws.Range("J15").Value = 35
$ws.Range("J18").Value = 53
$ws.Range("J19").Value = 95
$ws.Range("J20").Value = 65
$ws.Range("J27").Value = 14
$ws.Range("J29").Value = 149
$ws.Range("J33").Value = 130
$ws.Range("J37").Value = 109
$ws.Range("J42").Value = 145
$ws.Range("J54").Value = 55
$ws.Range("J63").Value = 17
$ws.Range("J65").Value = 77
$ws.Range("J67").Value = 116
$ws.Range("J68").Value = 9
$ws.Range("J69").Value = 6
$ws.Range("J70").Value = 8
$ws.Range("J76").Value = 51
$ws.Range("J78").Value = 40
$ws.Range("J79").Value = 98
$ws.Range("J83").Value = 59
$ws.Range("J84").Value = 35
$ws.Range("J88").Value = 23
$ws.Range("J89").Value = 38
$ws.Range("J90").Value = 38
$ws.Range("J94").Value = 19
$ws.Range("J95").Value = 57
$ws.Range("J99").Value = 40
$ws.Range("J101").Value = 3042

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J3").Value = 17
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 59

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J3").Value = 19
$ws.Range("J7").Value = 57

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 29
$ws.Range("J3").Value = 35
$ws.Range("J6").Value = 55
$ws.Range("J7").Value = 130

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J2").Value = 16
$ws.Range("J7").Value = 55

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 49
$ws.Range("J6").Value = 38
$ws.Range("J7").Value = 149

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J6").Value = 41
$ws.Range("J7").Value = 95

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J6").Value = 31
$ws.Range("J7").Value = 51

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("J2").Value = 6
$ws.Range("J7").Value = 28

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 22
$ws.Range("J3").Value = 21
$ws.Range("J4").Value = 5
$ws.Range("J6").Value = 94
$ws.Range("J7").Value = 145

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("J3").Value = 9
$ws.Range("J6").Value = 7
$ws.Range("J7").Value = 19

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J3").Value = 16
$ws.Range("J7").Value = 40

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("J6").Value = 4
$ws.Range("J7").Value = 6

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 28
$ws.Range("J4").Value = 7
$ws.Range("J6").Value = 29
$ws.Range("J7").Value = 98

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J3").Value = 22
$ws.Range("J6").Value = 20
$ws.Range("J7").Value = 65

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("J3").Value = 5
$ws.Range("J7").Value = 53

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 19

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J2").Value = 10
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J2").Value = 12
$ws.Range("J7").Value = 41

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("J3").Value = 7
$ws.Range("J7").Value = 22

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("J2").Value = 9
$ws.Range("J3").Value = 6
$ws.Range("J7").Value = 22

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("J2").Value = 4
$ws.Range("J7").Value = 8

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J3").Value = 6
$ws.Range("J4").Value = 1
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 23

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J3").Value = 63
$ws.Range("J6").Value = 61
$ws.Range("J7").Value = 191

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J6").Value = 8
$ws.Range("J7").Value = 14

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J6").Value = 19
$ws.Range("J7").Value = 38

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("J3").Value = 2
$ws.Range("J7").Value = 9

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J6").Value = 28
$ws.Range("J7").Value = 85

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("J6").Value = 7
$ws.Range("J7").Value = 16
